$d = $word.ActiveDocument

$find = $d.Content.Find
$find.Execute("功能整合（玩家与敌人，地图的交互），", $true, $false, $false, $false, $false, $true, 1, $false, "功能整合（玩家与敌人，地图的交互，buff测试），", 2)
